$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 701, shifting rows 701:787 down to 703:789
$ws.Rows("701:702").Insert()

# Fill new row 701 with its data
$ws.Cells.Item(701, 1).Value = 4
$ws.Cells.Item(701, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(701, 3).Value = "Los Lagos"
$ws.Cells.Item(701, 4).Value = 45142
$ws.Cells.Item(701, 5).Value = 10
$ws.Cells.Item(701, 6).Value = 100112006
$ws.Cells.Item(701, 7).Value = "Repollo"
$ws.Cells.Item(701, 8).Value = "Copenhague"
$ws.Cells.Item(701, 9).Value = "Primera"
$ws.Cells.Item(701, 10).Value = 500
$ws.Cells.Item(701, 11).Value = 1600
$ws.Cells.Item(701, 12).Value = 1600
$ws.Cells.Item(701, 13).Value = 1600
$ws.Cells.Item(701, 14).Value = "`$/unidad"
$ws.Cells.Item(701, 15).Value = "Región Metropolitana"
$ws.Cells.Item(701, 16).Value = 1600
$ws.Cells.Item(701, 17).Value = 1
$ws.Cells.Item(701, 18).Value = "Hortaliza"

# Fill new row 702 with its data
$ws.Cells.Item(702, 1).Value = 4
$ws.Cells.Item(702, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(702, 3).Value = "Los Lagos"
$ws.Cells.Item(702, 4).Value = 45142
$ws.Cells.Item(702, 5).Value = 10
$ws.Cells.Item(702, 6).Value = 100112006
$ws.Cells.Item(702, 7).Value = "Repollo"
$ws.Cells.Item(702, 8).Value = "Crespo record"
$ws.Cells.Item(702, 9).Value = "Primera"
$ws.Cells.Item(702, 10).Value = 1000
$ws.Cells.Item(702, 11).Value = 1500
$ws.Cells.Item(702, 12).Value = 1500
$ws.Cells.Item(702, 13).Value = 1500
$ws.Cells.Item(702, 14).Value = "`$/unidad"
$ws.Cells.Item(702, 15).Value = "Región Metropolitana"
$ws.Cells.Item(702, 16).Value = 1500
$ws.Cells.Item(702, 17).Value = 1
$ws.Cells.Item(702, 18).Value = "Hortaliza"
